$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the "New York" / "Montevideo" columns (Y and Z) ---
# Header row (row 1): swap the column titles
$yHeader = $ws.Range("Y1").Value()
$zHeader = $ws.Range("Z1").Value()
$ws.Range("Y1").Value = $zHeader
$ws.Range("Z1").Value = $yHeader

# Data rows 2-6: swap the Y/Z values for each row
for ($r = 2; $r -le 6; $r++) {
    $yVal = $ws.Range("Y$r").Value()
    $zVal = $ws.Range("Z$r").Value()
    $ws.Range("Y$r").Value = $zVal
    $ws.Range("Z$r").Value = $yVal
}

# --- AH6 was the last row's date cell (date-only format); now that a new
#     last row (7) is being added, AH6 reverts to the regular timestamp
#     format used by the other interior rows ---
$ws.Range("AH6").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Append new row 7 with the latest bunker price data ---
$ws.Range("A7").Value = 560
$ws.Range("B7").Value = 640
$ws.Range("C7").Value = 492
$ws.Range("D7").Value = 545
$ws.Range("E7").Value = 645
$ws.Range("F7").Value = 639
$ws.Range("G7").Value = 498
$ws.Range("H7").Value = 532
$ws.Range("I7").Value = 561
$ws.Range("J7").Value = 515
$ws.Range("K7").Value = 575
$ws.Range("L7").Value = 505
$ws.Range("M7").Value = 533
$ws.Range("N7").Value = 880
$ws.Range("O7").Value = 583
$ws.Range("P7").Value = 523
$ws.Range("Q7").Value = 498
$ws.Range("R7").Value = 537
$ws.Range("S7").Value = 582
$ws.Range("T7").Value = 653
$ws.Range("U7").Value = 589
$ws.Range("V7").Value = 485
$ws.Range("W7").Value = 545
$ws.Range("X7").Value = 530
$ws.Range("Y7").Value = 548
$ws.Range("Z7").Value = 531
$ws.Range("AA7").Value = 498
$ws.Range("AB7").Value = 542
$ws.Range("AC7").Value = 585.5
$ws.Range("AD7").Value = 508
$ws.Range("AE7").Value = 510
$ws.Range("AF7").Value = 529
$ws.Range("AG7").Value = 485
$ws.Range("AH7").Value = 45735
$ws.Range("AI7").Value = 502
$ws.Range("AJ7").Value = 546
$ws.Range("AK7").Value = 515
$ws.Range("AL7").Value = 742
$ws.Range("AM7").Value = 656
$ws.Range("AN7").Value = 617
$ws.Range("AO7").Value = 500
$ws.Range("AP7").Value = 628
$ws.Range("AQ7").Value = 752
$ws.Range("AR7").Value = 511
$ws.Range("AS7").Value = 492
$ws.Range("AT7").Value = 554
$ws.Range("AU7").Value = 569
$ws.Range("AV7").Value = 642

# New row's Date cell (AH7) keeps the date-only format (matching what AH6
# used to have before it became an interior row)
$ws.Range("AH7").NumberFormat = "YYYY-MM-DD"
